{"js": "const replacements = [\n  [\"2025-07-17 Thursday\", \"2025-07-18 Friday\"],\n  [\"397\u00d74=\", \"362\u00d79=\"],\n  [\"945\u00d73=\", \"250\u00d78=\"],\n  [\"931\u00d79=\", \"639\u00d77=\"],\n  [\"148\u00d73=\", \"364\u00d76=\"],\n  [\"963\u00d76=\", \"815\u00d78=\"],\n  [\"672\u00d78=\", \"823\u00d78=\"],\n  [\"907\u00d73=\", \"326\u00d78=\"],\n  [\"556\u00d75=\", \"297\u00d74=\"],\n  [\"317\u00d72=\", \"722\u00d78=\"],\n  [\"557\u00d76=\", \"788\u00d79=\"],\n  [\"792\u00d74=\", \"754\u00d73=\"],\n  [\"649\u00d73=\", \"120\u00d76=\"],\n  [\"498\u00d78=\", \"683\u00d77=\"],\n  [\"308\u00d75=\", \"860\u00d73=\"],\n  [\"459\u00d77=\", \"985\u00d72=\"],\n  [\"749\u00d77=\", \"169\u00d77=\"],\n  [\"302\u00d74=\", \"395\u00d78=\"],\n  [\"337\u00d78=\", \"765\u00d75=\"],\n  [\"248\u00d78=\", \"607\u00d73=\"],\n  [\"185\u00d77=\", \"784\u00d77=\"],\n  [\"338\u00d74=\", \"180\u00d79=\"],\n  [\"342\u00d77=\", \"963\u00d77=\"],\n  [\"356\u00d76=\", \"765\u00d72=\"],\n  [\"382\u00d77=\", \"464\u00d73=\"],\n  [\"829\u00d76=\", \"479\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-17 Thursday\", \"2025-07-18 Friday\"),\n    @(\"397\u00d74=\", \"362\u00d79=\"),\n    @(\"945\u00d73=\", \"250\u00d78=\"),\n    @(\"931\u00d79=\", \"639\u00d77=\"),\n    @(\"148\u00d73=\", \"364\u00d76=\"),\n    @(\"963\u00d76=\", \"815\u00d78=\"),\n    @(\"672\u00d78=\", \"823\u00d78=\"),\n    @(\"907\u00d73=\", \"326\u00d78=\"),\n    @(\"556\u00d75=\", \"297\u00d74=\"),\n    @(\"317\u00d72=\", \"722\u00d78=\"),\n    @(\"557\u00d76=\", \"788\u00d79=\"),\n    @(\"792\u00d74=\", \"754\u00d73=\"),\n    @(\"649\u00d73=\", \"120\u00d76=\"),\n    @(\"498\u00d78=\", \"683\u00d77=\"),\n    @(\"308\u00d75=\", \"860\u00d73=\"),\n    @(\"459\u00d77=\", \"985\u00d72=\"),\n    @(\"749\u00d77=\", \"169\u00d77=\"),\n    @(\"302\u00d74=\", \"395\u00d78=\"),\n    @(\"337\u00d78=\", \"765\u00d75=\"),\n    @(\"248\u00d78=\", \"607\u00d73=\"),\n    @(\"185\u00d77=\", \"784\u00d77=\"),\n    @(\"338\u00d74=\", \"180\u00d79=\"),\n    @(\"342\u00d77=\", \"963\u00d77=\"),\n    @(\"356\u00d76=\", \"765\u00d72=\"),\n    @(\"382\u00d77=\", \"464\u00d73=\"),\n    @(\"829\u00d76=\", \"479\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
